# Add 2022-Q1 data:
#  - rename current "总计" sheet to "2022-Q1" and fill it with the fund
#    detail rows for the new quarter (same layout as 2021-Q1/Q2/Q3)
#  - append a brand-new "总计" sheet at the end with the summary table,
#    now including a 2022-Q1 row on top of the previous quarters

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Turn the existing "总计" sheet into the new "2022-Q1" detail sheet
# ---------------------------------------------------------------------
$fundSheet = $wb.Worksheets.Item("总计")
$fundSheet.Name = "2022-Q1"
$fundSheet.Cells.Clear()

# Header row
$fundSheet.Range("B1").Value = "基金代码"
$fundSheet.Range("C1").Value = "基金名称"
$fundSheet.Range("D1").Value = "基金规模"
$fundSheet.Range("E1").Value = "股票总仓位"
$fundSheet.Range("F1").Value = "仓位占比"
$fundSheet.Range("G1").Value = "持有市值(亿元)"
$fundSheet.Range("H1").Value = "仓位排名"

$headerRange = $fundSheet.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Index column (A) style
$indexRange = $fundSheet.Range("A2:A3")
$indexRange.Font.Bold = $true
$indexRange.HorizontalAlignment = -4108
$indexRange.VerticalAlignment = -4160
$indexRange.Borders.LineStyle = 1

$fundSheet.Range("A2").Value = 0
$fundSheet.Range("A3").Value = 1

# Numeric-looking text columns (fund code / amounts) must stay text so
# leading zeros and exact decimal strings are preserved, not converted
# to numbers. Temporarily force a text format, assign, then restore the
# number format back to General so no stray formatting is left behind.
$numericTextCols = "B", "D", "E", "F", "G"
foreach ($col in $numericTextCols) {
    $fundSheet.Range($col + "2:" + $col + "3").NumberFormat = "@"
}

$fundSheet.Range("B2").Value = "000369"
$fundSheet.Range("C2").Value = "广发全球医疗保健(QDII) - 人民币"
$fundSheet.Range("D2").Value = "2.46"
$fundSheet.Range("E2").Value = "81.85"
$fundSheet.Range("F2").Value = "3.48"
$fundSheet.Range("G2").Value = "0.0856"
$fundSheet.Range("H2").Value = 4

$fundSheet.Range("B3").Value = "000370"
$fundSheet.Range("C3").Value = "广发全球医疗保健(QDII) - 美元"
$fundSheet.Range("D3").Value = "2.46"
$fundSheet.Range("E3").Value = "81.85"
$fundSheet.Range("F3").Value = "3.48"
$fundSheet.Range("G3").Value = "0.0856"
$fundSheet.Range("H3").Value = 4

foreach ($col in $numericTextCols) {
    $fundSheet.Range($col + "2:" + $col + "3").NumberFormat = "General"
}

# ---------------------------------------------------------------------
# 2. Add a fresh "总计" sheet at the end with the updated summary table
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$totalSheet.Name = "总计"

$totalSheet.Range("B1").Value = "日期"
$totalSheet.Range("C1").Value = "持有数量(只)"
$totalSheet.Range("D1").Value = "持有市值(亿元)"

$totalHeaderRange = $totalSheet.Range("B1:D1")
$totalHeaderRange.Font.Bold = $true
$totalHeaderRange.HorizontalAlignment = -4108
$totalHeaderRange.VerticalAlignment = -4160
$totalHeaderRange.Borders.LineStyle = 1

$totalIndexRange = $totalSheet.Range("A2:A5")
$totalIndexRange.Font.Bold = $true
$totalIndexRange.HorizontalAlignment = -4108
$totalIndexRange.VerticalAlignment = -4160
$totalIndexRange.Borders.LineStyle = 1

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.17

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2021-Q3"
$totalSheet.Range("C3").Value = 4
$totalSheet.Range("D3").Value = 0.13

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2021-Q2"
$totalSheet.Range("C4").Value = 2
$totalSheet.Range("D4").Value = 0.11

$totalSheet.Range("A5").Value = 3
$totalSheet.Range("B5").Value = "2021-Q1"
$totalSheet.Range("C5").Value = 2
$totalSheet.Range("D5").Value = 0.1

# Keep the originally active/selected tab ("2021-Q1") active, since adding
# worksheets shifts the active tab to the newly created sheet by default.
$wb.Worksheets.Item(1).Activate()
